$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '88.412.73'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.277.98'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.34'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.55%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '631.23'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.76%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.390'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +17.78%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.720'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +16.59%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.275.73'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.583'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.49%  '
$ws.Range('E11').ClearFormats()
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'TRON'
$ws.Range('B12').ClearFormats()
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C12').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.187'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +11.78%  '
$ws.Range('E12').ClearFormats()
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('B13').ClearFormats()
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('C13').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000270'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.10%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.67'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.53'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.881.83'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.171.31'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.290.27'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.30'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.66%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.17'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.48%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '440.55'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -6.87%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.05'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.38'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.94%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.41'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.36'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.50'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -7.35%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000144'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +11.84%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.451.81'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '77.90'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.181'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -17.13%  '
$ws.Range('E31').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('E32').ClearFormats()
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('B33').ClearFormats()
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('C33').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '576.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.82%  '
$ws.Range('E33').ClearFormats()
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('B34').ClearFormats()
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C34').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '8.94'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.92%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.44'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +9.79%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.39'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -10.08%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.98'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.79%  '
$ws.Range('E37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.141'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.89%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.17'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.25%  '
$ws.Range('E39').ClearFormats()
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('B40').ClearFormats()
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C40').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.33'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.93%  '
$ws.Range('E40').ClearFormats()
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'WhiteBITCoin'
$ws.Range('B41').ClearFormats()
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('C41').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '21.84'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.62%  '
$ws.Range('E41').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.407'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.39%  '
$ws.Range('E43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.06'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.64%  '
$ws.Range('E44').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.72'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.11%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.140'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +23.34%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '181.34'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.65%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.36'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.22'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.57%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.29'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.63%  '
$ws.Range('E51').ClearFormats()
